$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws4 = $wb.Worksheets.Item(4)

# Row 3
$ws1.Cells.Item(2, 1).Copy() | Out-Null
$ws1.Cells.Item(3, 1).PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(3, 1).Value = 2
$ws1.Range("B3").Value = "'2024-11-09"
$ws1.Range("C3").Value = "安徽·崩坏同人only 爱莉希雅同人生日会"
$ws1.Range("D3").Value = "徽州大道与扬子江路口天琅百老汇一楼123号 禧棠捌号XITANGBH禧宴中心"
$ws1.Range("E3").Value = "'2024.11.09 12:00-11.09 22:00"
$ws1.Range("F3").Value = 5
$ws1.Range("G3").Value = 49
$ws1.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=93461"
$ws1.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202410/VnEQZYTQ1728892328769.png"
$ws1.Range("B3:I3").ClearFormats()

# Row 4
$ws1.Cells.Item(3, 1).Copy() | Out-Null
$ws1.Cells.Item(4, 1).PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(4, 1).Value = 3
$ws1.Range("B4").Value = "'2024-11-16"
$ws1.Range("C4").Value = "合肥·第九届环形宇宙动漫游戏嘉年华"
$ws1.Range("D4").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws1.Range("E4").Value = "'2024.11.16 09:30-11.17 17:00"
$ws1.Range("F4").Value = 3265
$ws1.Range("G4").Value = 72
$ws1.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=92565"
$ws1.Range("I4").Value = "//i1.hdslb.com/bfs/openplatform/202410/sxfiGFBq1728715876124.jpeg"
$ws1.Range("B4:I4").ClearFormats()

# Row 5
$ws1.Cells.Item(4, 1).Copy() | Out-Null
$ws1.Cells.Item(5, 1).PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(5, 1).Value = 4
$ws1.Range("B5").Value = "'2024-11-17"
$ws1.Range("C5").Value = "合肥·MAX特摄同人only2.0"
$ws1.Range("D5").Value = "桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间"
$ws1.Range("E5").Value = "'2024.11.17 10:00-11.17 18:00"
$ws1.Range("F5").Value = 153
$ws1.Range("G5").Value = 60
$ws1.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=92498"
$ws1.Range("I5").Value = "//i1.hdslb.com/bfs/openplatform/202410/ccEfc1521728888008037.jpeg"
$ws1.Range("B5:I5").ClearFormats()

# Row 6
$ws1.Cells.Item(5, 1).Copy() | Out-Null
$ws1.Cells.Item(6, 1).PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(6, 1).Value = 5
$ws1.Range("B6").Value = "'2024-11-30"
$ws1.Range("C6").Value = "合肥·风月引代号鸢同人only"
$ws1.Range("D6").Value = "徽州大道与杨子江路交口天琅百老汇一楼123号 禧棠捌号禧宴中心（滨湖店）"
$ws1.Range("E6").Value = "'2024.11.30 10:00-11.30 21:00"
$ws1.Range("F6").Value = 27
$ws1.Range("G6").Value = 55
$ws1.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=93322"
$ws1.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202409/Tu5YLbGx1727179854562.jpeg"
$ws1.Range("B6:I6").ClearFormats()

# Row 7
$ws1.Cells.Item(6, 1).Copy() | Out-Null
$ws1.Cells.Item(7, 1).PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(7, 1).Value = 6
$ws1.Range("B7").Value = "'2024-12-07"
$ws1.Range("C7").Value = "合肥·心动恋章·冬日序国乙&代号鸢同人only"
$ws1.Range("D7").Value = "上海路与迎淮路交口向东200米(云峯中心一楼) 费加罗宴会艺术中心(省府店)"
$ws1.Range("E7").Value = "'2024.12.07 12:00-12.07 21:00"
$ws1.Range("F7").Value = 144
$ws1.Range("G7").Value = 50
$ws1.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=93319"
$ws1.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202409/KtMLL8ZO1727684987784.jpeg"
$ws1.Range("B7:I7").ClearFormats()

# Row 7
$ws4.Cells.Item(6, 1).Copy() | Out-Null
$ws4.Cells.Item(7, 1).PasteSpecial(-4122) | Out-Null
$ws4.Cells.Item(7, 1).Value = 6
$ws4.Range("B7").Value = "'2024-11-09"
$ws4.Range("C7").Value = "安徽·崩坏同人only 爱莉希雅同人生日会"
$ws4.Range("D7").Value = "徽州大道与扬子江路口天琅百老汇一楼123号 禧棠捌号XITANGBH禧宴中心"
$ws4.Range("E7").Value = "'2024.11.09 12:00-11.09 22:00"
$ws4.Range("F7").Value = 5
$ws4.Range("G7").Value = 49
$ws4.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=93461"
$ws4.Range("I7").Value = "//i2.hdslb.com/bfs/openplatform/202410/VnEQZYTQ1728892328769.png"
$ws4.Range("B7:I7").ClearFormats()

# Row 8
$ws4.Cells.Item(7, 1).Copy() | Out-Null
$ws4.Cells.Item(8, 1).PasteSpecial(-4122) | Out-Null
$ws4.Cells.Item(8, 1).Value = 7
$ws4.Range("B8").Value = "'2024-11-16"
$ws4.Range("C8").Value = "合肥·第九届环形宇宙动漫游戏嘉年华"
$ws4.Range("D8").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws4.Range("E8").Value = "'2024.11.16 09:30-11.17 17:00"
$ws4.Range("F8").Value = 3265
$ws4.Range("G8").Value = 72
$ws4.Range("H8").Value = "https://show.bilibili.com/platform/detail.html?id=92565"
$ws4.Range("I8").Value = "//i1.hdslb.com/bfs/openplatform/202410/sxfiGFBq1728715876124.jpeg"
$ws4.Range("B8:I8").ClearFormats()

# Row 9
$ws4.Cells.Item(8, 1).Copy() | Out-Null
$ws4.Cells.Item(9, 1).PasteSpecial(-4122) | Out-Null
$ws4.Cells.Item(9, 1).Value = 8
$ws4.Range("B9").Value = "'2024-11-17"
$ws4.Range("C9").Value = "合肥·MAX特摄同人only2.0"
$ws4.Range("D9").Value = "桐城路127号合作经济广场3号楼23层 赤阑桥艺术空间"
$ws4.Range("E9").Value = "'2024.11.17 10:00-11.17 18:00"
$ws4.Range("F9").Value = 153
$ws4.Range("G9").Value = 60
$ws4.Range("H9").Value = "https://show.bilibili.com/platform/detail.html?id=92498"
$ws4.Range("I9").Value = "//i1.hdslb.com/bfs/openplatform/202410/ccEfc1521728888008037.jpeg"
$ws4.Range("B9:I9").ClearFormats()

# Row 10
$ws4.Cells.Item(9, 1).Copy() | Out-Null
$ws4.Cells.Item(10, 1).PasteSpecial(-4122) | Out-Null
$ws4.Cells.Item(10, 1).Value = 9
$ws4.Range("B10").Value = "'2024-11-30"
$ws4.Range("C10").Value = "合肥·风月引代号鸢同人only"
$ws4.Range("D10").Value = "徽州大道与杨子江路交口天琅百老汇一楼123号 禧棠捌号禧宴中心（滨湖店）"
$ws4.Range("E10").Value = "'2024.11.30 10:00-11.30 21:00"
$ws4.Range("F10").Value = 27
$ws4.Range("G10").Value = 55
$ws4.Range("H10").Value = "https://show.bilibili.com/platform/detail.html?id=93322"
$ws4.Range("I10").Value = "//i1.hdslb.com/bfs/openplatform/202409/Tu5YLbGx1727179854562.jpeg"
$ws4.Range("B10:I10").ClearFormats()

# Row 11
$ws4.Cells.Item(10, 1).Copy() | Out-Null
$ws4.Cells.Item(11, 1).PasteSpecial(-4122) | Out-Null
$ws4.Cells.Item(11, 1).Value = 10
$ws4.Range("B11").Value = "'2024-12-07"
$ws4.Range("C11").Value = "合肥·一生必听的古典系列《钟》—超技钢琴曲炫彩音乐会"
$ws4.Range("D11").Value = "徽州大道辅路与祁门路辅路交叉口北120米 包河凤凰剧院"
$ws4.Range("E11").Value = "'2024.12.07 19:30-12.07 21:00"
$ws4.Range("F11").Value = 2
$ws4.Range("G11").Value = 72
$ws4.Range("H11").Value = "https://show.bilibili.com/platform/detail.html?id=91608"
$ws4.Range("I11").Value = "//i0.hdslb.com/bfs/openplatform/202408/wiLiWoeM1725005636569.jpeg"
$ws4.Range("B11:I11").ClearFormats()

# Row 12
$ws4.Cells.Item(11, 1).Copy() | Out-Null
$ws4.Cells.Item(12, 1).PasteSpecial(-4122) | Out-Null
$ws4.Cells.Item(12, 1).Value = 11
$ws4.Range("B12").Value = "'2024-12-07"
$ws4.Range("C12").Value = "合肥·心动恋章·冬日序国乙&代号鸢同人only"
$ws4.Range("D12").Value = "上海路与迎淮路交口向东200米(云峯中心一楼) 费加罗宴会艺术中心(省府店)"
$ws4.Range("E12").Value = "'2024.12.07 12:00-12.07 21:00"
$ws4.Range("F12").Value = 144
$ws4.Range("G12").Value = 50
$ws4.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=93319"
$ws4.Range("I12").Value = "//i0.hdslb.com/bfs/openplatform/202409/KtMLL8ZO1727684987784.jpeg"
$ws4.Range("B12:I12").ClearFormats()
